$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (changed) date column (C) for rows 2-6
# from serial date 45174 (2023-09-05) to 45175 (2023-09-06)
foreach ($r in 2..6) {
    $ws.Cells.Item($r, 3).Value = 45175
}
